# Update the "dSF" (column F) values on Sheet1 to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -2
    8  = 3
    11 = 0
    16 = 0
    20 = 4
    21 = 2
    33 = -1
    37 = -2
    39 = 3
    40 = -1
    41 = 0
    44 = 7
    45 = -3
    47 = 0
    49 = -3
    54 = -5
    55 = 4
    58 = -6
    61 = -9
    65 = -8
    66 = -2
    69 = -1
    71 = 6
    73 = -4
    74 = 0
    78 = -2
    80 = -2
    83 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
